$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$sec = $d.Sections(1)

# --- Footers: both the primary footer and the first-page footer contain the
# Pearson Edexcel logo picture whose internal part name is "image2.png";
# rename it to "image1.png" (docPr/name + cNvPr/name in the underlying XML).
$primaryFooterShapes = $sec.Footers(1).Range.InlineShapes
for ($i = 1; $i -le $primaryFooterShapes.Count; $i++) {
    $shp = $primaryFooterShapes.Item($i)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image1.png"
    }
}

$firstFooterShapes = $sec.Footers(2).Range.InlineShapes
for ($i = 1; $i -le $firstFooterShapes.Count; $i++) {
    $shp = $firstFooterShapes.Item($i)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image1.png"
    }
}

# --- Header: the first-page header contains the BTEC logo picture whose
# internal part name is "image1.jpg"; rename it to "image2.jpg".
$firstHeaderShapes = $sec.Headers(2).Range.InlineShapes
for ($i = 1; $i -le $firstHeaderShapes.Count; $i++) {
    $shp = $firstHeaderShapes.Item($i)
    if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
        $shp.Name = "image2.jpg"
    }
}
